$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove all existing rows (clears stale cell content, styles-on-empty-cells, and row heights)
# so the sheet can be rebuilt to match the new shape exactly (24 rows -> 23 rows).
$ws.Rows("1:24").Delete()

# Row 1
$ws.Range('B1').Value = 'Ementa atual:'
$ws.Range('C1').Value = 'Ementa modificada (dados modificados em vermelho):'

# Row 2
$ws.Range('B2').Value = 'LOQ4050'
$ws.Range('C2').Value = 'LOQ4050'

# Row 3
$ws.Range('A3').Value = 'Nome:'
$ws.Range('B3').Value = ' Engenharia Econômica'
$ws.Range('C3').Value = ' Engenharia Econômica'

# Row 4
$ws.Range('A4').Value = 'Name:'
$ws.Range('B4').Value = 'Economic Engineering'
$ws.Range('C4').Value = 'Economic Engineering'

# Row 5
$ws.Range('A5').Value = 'Créditos-aula:'
$ws.Range('B5').Value = '2'
$ws.Range('C5').Value = '2'

# Row 6
$ws.Range('A6').Value = 'Créditos-trabalho'
$ws.Range('B6').Value = '0'
$ws.Range('C6').Value = '0'

# Row 7
$ws.Range('A7').Value = 'Carga horária:'
$ws.Range('B7').Value = '30 h'
$ws.Range('C7').Value = '30 h'

# Row 8
$ws.Range('A8').Value = 'Ativação:'
$ws.Range('B8').Value = '01/01/2011'
$ws.Range('C8').Value = '01/01/2011'

# Row 9
$ws.Range('A9').Value = 'Semestre ideal:'
$ws.Range('B9').Value = 'EF-7,EQD-9,EQN-10'
$ws.Range('C9').Value = 'EF-7,EQD-9,EQN-10'

# Row 10
$ws.Range('A10').Value = 'Objetivos:'
$ws.Range('B10').Value = '5840671 - Francisco José Moreira Chaves'
$ws.Range('C10').Value = '5840671 - Francisco José Moreira Chaves'
$ws.Rows.Item(10).RowHeight = 60

# Row 11
$ws.Range('A11').Value = 'Objectives:'
$ws.Range('B11').Value = '1 ) Formative : Offer the learner the basic conditions necessary for their vocational training.2 ) Informational : provide the student the basics to understanding, advice and monitoring of Projects in the Chemical Industry following specific methodology .3 ) :The learner develop analytical reasoning , following systematic methodology applied in projects .'
$ws.Range('C11').Value = '1 ) Formative : Offer the learner the basic conditions necessary for their vocational training.2 ) Informational : provide the student the basics to understanding, advice and monitoring of Projects in the Chemical Industry following specific methodology .3 ) :The learner develop analytical reasoning , following systematic methodology applied in projects .'
$ws.Rows.Item(11).RowHeight = 60

# Row 12
$ws.Range('A12').Value = 'Docentes responsáveis:'

# Row 13
$ws.Range('A13').Value = 'Programa resumido:'
$ws.Range('B13').Value = 'Semestral'
$ws.Range('C13').Value = 'Semestral'
$ws.Rows.Item(13).RowHeight = 60

# Row 14
$ws.Range('A14').Value = 'Short syllabus:'
$ws.Range('B14').Value = 'Market - Estimated investment : - Economic Analysis of Investments'
$ws.Range('C14').Value = 'Market - Estimated investment : - Economic Analysis of Investments'
$ws.Rows.Item(14).RowHeight = 60

# Row 15
$ws.Range('A15').Value = 'Programa:'
$ws.Range('B15').Value = '01/01/2011'
$ws.Range('C15').Value = '01/01/2011'
$ws.Rows.Item(15).RowHeight = 120

# Row 16
$ws.Range('A16').Value = 'Syllabus:'
$ws.Range('B16').Value = 'Market - market trends - Marketing and market analysis - Estimate of investment: - working capital - capital assets - tangible / intangibles assets - Preliminary investments / real investments - fixed costs / variable costs , depreciation , economic concept of externalities and approaches theoretical , elements to internalize externalities , direct control or stabilize the company , indirect methods and observed data , indirect methods and data assumptions , methods and alleged direct data  methods and observed data , Economic Analysis of Investments ( basics aspects) - Statistical Techniques and interest rates - basics of Engineering Economy - Beta Distribution for analysis in the risk environment .'
$ws.Range('C16').Value = 'Market - market trends - Marketing and market analysis - Estimate of investment: - working capital - capital assets - tangible / intangibles assets - Preliminary investments / real investments - fixed costs / variable costs , depreciation , economic concept of externalities and approaches theoretical , elements to internalize externalities , direct control or stabilize the company , indirect methods and observed data , indirect methods and data assumptions , methods and alleged direct data  methods and observed data , Economic Analysis of Investments ( basics aspects) - Statistical Techniques and interest rates - basics of Engineering Economy - Beta Distribution for analysis in the risk environment .'
$ws.Rows.Item(16).RowHeight = 120

# Row 17
$ws.Range('A17').Value = 'Avaliação:'

# Row 18
$ws.Range('A18').Value = 'Método:'
$ws.Range('B18').Value = '5840671 - Francisco José Moreira Chaves'
$ws.Range('C18').Value = '5840671 - Francisco José Moreira Chaves'
$ws.Rows.Item(18).RowHeight = 60

# Row 19
$ws.Range('A19').Value = 'Critério:'
$ws.Range('B19').Value = 'Por meio de aulas presenciais, com apresentação dos fundamentos, e resolução de exercícios e exemplos aplicativos com uso de tabelas e normas específicas.'
$ws.Range('C19').Value = 'Por meio de aulas presenciais, com apresentação dos fundamentos, e resolução de exercícios e exemplos aplicativos com uso de tabelas e normas específicas.'
$ws.Rows.Item(19).RowHeight = 60

# Row 20
$ws.Range('A20').Value = 'Norma de recuperação:'
$ws.Range('B20').Value = 'A Avaliação será: (P1 + 2P2)/3 Onde: P1: Prova Individual - c/ peso-1 P2: : Prova Individual - c/ peso-2'
$ws.Range('C20').Value = 'A Avaliação será: (P1 + 2P2)/3 Onde: P1: Prova Individual - c/ peso-1 P2: : Prova Individual - c/ peso-2'
$ws.Rows.Item(20).RowHeight = 60

# Row 21
$ws.Range('A21').Value = 'Bibliografia:'
$ws.Range('B21').Value = 'Prova de exame.'
$ws.Range('C21').Value = 'Prova de exame.'
$ws.Rows.Item(21).RowHeight = 120

# Row 22
$ws.Range('A22').Value = 'Requisitos:'

# Row 23
$ws.Range('B23').Value = 'LOB1012 -  Estatística  (Requisito)
'
$ws.Range('C23').Value = 'LOB1012 -  Estatística  (Requisito)
'
$ws.Rows.Item(23).RowHeight = 30

